$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was added to the "Puerro" price sheet. It belongs
# right above the current row 40 (chronologically it is the most recent
# record for this series), so insert a blank row there which pushes the
# existing data (old rows 40-118) down to rows 41-119, then populate the
# newly inserted row with its own values.
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value2 = 9
$ws.Range("B40").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value2 = "Metropolitana"
$ws.Range("D40").Value2 = 45028
$ws.Range("E40").Value2 = 13
$ws.Range("F40").Value2 = 100112005
$ws.Range("G40").Value2 = "Puerro"
$ws.Range("H40").Value2 = "Sin especificar"
$ws.Range("I40").Value2 = "Primera"
$ws.Range("J40").Value2 = 160
$ws.Range("K40").Value2 = 8000
$ws.Range("L40").Value2 = 8000
$ws.Range("M40").Value2 = 8000
$ws.Range("N40").Value2 = "`$/paquete 20 unidades"
$ws.Range("O40").Value2 = "Provincia de Chacabuco"
$ws.Range("P40").Value2 = 400
$ws.Range("Q40").Value2 = 20
$ws.Range("R40").Value2 = "Hortaliza"
